$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window / view cosmetics -------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 88

# --- Row 3 : ILS columns (J/K/L) + matrix-time columns (D/G) ---------------
$ws.Range("D3").Value = 0.0002
$ws.Range("G3").Value = 0.0013
$ws.Range("J3").Value = 2.5406
$ws.Range("K3").Value = 27603
$ws.Range("L3").Formula = "=(K3-C3)/K3"

# --- Row 4 -------------------------------------------------------------------
$ws.Range("D4").Value = 0.0008
$ws.Range("G4").Value = 0.0031
$ws.Range("J4").NumberFormat = "#,##0.00"
$ws.Range("J4").Value = 4.8602
$ws.Range("K4").Value = 6808
$ws.Range("L4").Formula = "=(K4-C4)/K4"

# --- Row 5 -------------------------------------------------------------------
$ws.Range("D5").Value = 0.0043
$ws.Range("G5").Value = 0.1804
$ws.Range("J5").Value = 32.4357
$ws.Range("K5").Value = 10830
$ws.Range("L5").Formula = "=(K5-C5)/K5"

# --- Row 6 -------------------------------------------------------------------
$ws.Range("D6").Value = 0.0531
$ws.Range("G6").Value = 5.6136
$ws.Range("J6").Value = 79.8824
$ws.Range("K6").Value = 95007
$ws.Range("L6").Formula = "=(K6-C6)/K6"

# --- Row 7 : only the matrix-time column changes ------------------------------
$ws.Range("D7").Value = 4.3745

# --- Row 8 ---------------------------------------------------------------------
$ws.Range("D8").Value = 9.2075

# --- Row 9 ---------------------------------------------------------------------
$ws.Range("D9").Value = 11.552

# --- Comments on D7 / D8 ------------------------------------------------------
$excel.UserName = "JuanG"
$note = "JuanG:" + [char]10 + "No toma en cuenta el tiempo que tarda en generar la matriz de distancias"
$ws.Range("D7").AddComment($note)
$ws.Range("D8").AddComment($note)
